# BOM.xlsx update: add two new jumper rows (JP401/JP402) to the "Passives"
# sheet, and restore the on-screen selection/scroll state that Excel wrote
# out when the workbook was last saved by the author.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ICs")
$ws2 = $wb.Worksheets.Item("Passives")

# --- New data rows on the "Passives" sheet ---------------------------------
# Columns: A=Component Number, B=Purpose (Hypothesis), C=Description,
#          D=Type, E=Voltage (V) (After) [Resistors only],
#          F=Connection 1, G=Connection 2
#
# Cells are written in this specific order so that newly-created shared
# strings land at the same indices the workbook author's Excel produced
# (JP401, LH7A400 T1, LH7A400 T9, JP402, then the repeated "Jumper" last).
$ws2.Range("A9").Value = "JP401"
$ws2.Range("F9").Value = "LH7A400 T1"
$ws2.Range("G9").Value = "LH7A400 T9"
$ws2.Range("A10").Value = "JP402"
$ws2.Range("D9").Value = "Jumper"
$ws2.Range("D10").Value = "Jumper"

# --- Restore view state ------------------------------------------------
# Passives sheet: remembered selection moves from B2 to D10.
[void]$ws2.Range("D10").Select()

# ICs sheet stays the active tab; its remembered selection moves from A3
# to J10, scrolled so row 8 is at the top of the viewport.
[void]$ws1.Activate()
[void]$ws1.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("J10").Select()

# Workbook window position/size, as last saved by the author.
$excel.ActiveWindow.Left = 220
$excel.ActiveWindow.Top = 460
$excel.ActiveWindow.Width = 30960
$excel.ActiveWindow.Height = 20540
